# Daily attendance processing - 2025-12-04 13:43:20
#
# Normalizes the "Recorded By" column (G) so that system-generated
# recorder tags ("system" / "System") are listed before any human/user
# identifiers (emails) in the comma-separated list, while leaving rows
# that were recorded by an administrator (admin@admin.com) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $text = [string]$current
    if ($text -eq "") {
        continue
    }

    # Leave administrator-recorded rows untouched.
    if ($text.Contains("admin@admin.com")) {
        continue
    }

    $parts = $text.Split(",")
    if ($parts.Count -lt 2) {
        continue
    }

    $systemTokens = @()
    $otherTokens = @()

    foreach ($part in $parts) {
        $trimmed = $part.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemTokens += $trimmed
        } else {
            $otherTokens += $trimmed
        }
    }

    if ($systemTokens.Count -eq 0) {
        continue
    }

    $ordered = $systemTokens + $otherTokens
    $newText = $ordered -join ", "

    if ($newText -ne $text) {
        $cell.Value2 = $newText
    }
}
